# The source sheet had a duplicated "extraTurn" row (row 16) under the
# Player attribute block - it duplicated row 8's extraTurn entry but with
# a wrong offset (42 instead of 6). Remove that duplicate row; Excel will
# shift every row below it up by one and fix up the row-relative formulas
# (shared formula ranges, the relative K43 = J43 formula, dimension, etc.)
# automatically, exactly as described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")
$ws.Activate()

$ws.Rows.Item(16).Delete()

# Restore the view/selection state expected after the edit.
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("G12:K12").Select()

$wb.Save()
